$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1) re-styling.
#    A1 -> style of A2 (size-14 font), B1 -> style of B2 (size-16 font),
#    C1 keeps its current style, D1:I1 -> style of D2 (size-16 font).
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B1").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("D1:I1").PasteSpecial(-4122)

$ws.Rows.Item(1).RowHeight = 46.8

# ---------------------------------------------------------------------------
# 2. Drop now-unused ACTUAL RESULT / STATUS cells (H,I) and the blank
#    TEST DATA cells (F) for the existing test cases (rows 2-8).
# ---------------------------------------------------------------------------
$ws.Range("F2:F8").Clear()
$ws.Range("H2:I8").Clear()

# ---------------------------------------------------------------------------
# 3. Three new test cases (rows 9-11).
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A9:A11").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B9:B11").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("D9:D11").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("E9:E11").PasteSpecial(-4122)

$ws.Range("G2").Copy()
$ws.Range("G9:G11").PasteSpecial(-4122)

# Fields are entered in the same left-to-right, top-to-bottom order the
# author used, since that order controls how new entries land in the
# shared-string table (matching cell text is de-duplicated/reused).
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Verify error message appear for non-existing email"
$ws.Range("D9").Value = "User on password recovery page"
$ws.Range("E9").Value = "1-Navigate to login page`n2-Click on forgot password link`n3-Enter non-existing email`n4-Click on ""Password recovery"" button"
$ws.Range("G9").Value = "1-An error message shall appear`n2-No email will be sent"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Verify error message for expired password reset link"
$ws.Range("D10").Value = "User on password recovery page"
$ws.Range("E10").Value = "1-Navigate to login page`n2-Click on forgot `npassword link`n3-Wait for 24 Hour`n4-Click on ""Password recovery"" button"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Verify error message for inequality new password and it's confirmation"
$ws.Range("D11").Value = "User on password recovery page"
$ws.Range("E11").Value = "1-Enter a new password `n2-Enter a different password at the ""confirm password"" field"

$ws.Range("G10").Value = "An error message shall appear"
$ws.Range("G11").Value = "An error message shall appear"

$ws.Rows.Item(9).RowHeight = 189
$ws.Rows.Item(10).RowHeight = 168
$ws.Rows.Item(11).RowHeight = 126

# ---------------------------------------------------------------------------
# 4. Column C got wider now that it holds its own test-data notes.
# ---------------------------------------------------------------------------
$ws.Range("C:C").ColumnWidth = 19.2

# ---------------------------------------------------------------------------
# 5. View state: a bit of zoom-out and the selection sitting below the
#    freshly added rows.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 96
$ws.Range("A12").Select()
